# New logout test case: LOGOUT-02 (restricted page access after logout)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logout Test")

# Make this sheet the active one (it was already tabSelected in the source file)
$ws.Activate()

# Copy the existing LOGOUT-01 row (row 4) down to row 5 so the new row
# inherits the same styles/formatting (fonts, fills, borders, alignment),
# then overwrite the values for the new test case.
$ws.Range("C4:M4").Copy($ws.Range("C5:M5"))
$excel.CutCopyMode = $false

# Fill in the ID-like columns first, then the rest, matching the order the
# cells were authored in (SUB, TC ID, Scenario Ref, Scenario, Test
# Description, Precondition, Steps, Test Data); Module/Expected/Priority
# reuse text already used by LOGOUT-01.
$ws.Range("C5").Value = "LOGOUT-02"
$ws.Range("E5").Value = "LOGOUT-TC-02"
$ws.Range("G5").Value = "LOGOUT-TS-02"
$ws.Range("D5").Value = "Verify user cannot access restricted page after logout"
$ws.Range("H5").Value = "Access resticted page after logout"
$ws.Range("I5").Value = "User already logged out"
$ws.Range("J5").Value = "1. Logout `n2. Enter restricted page URL in browser `n3. Press Enter"
$ws.Range("K5").Value = "Directed restricted URL"
$ws.Range("F5").Value = "Logout Module"
$ws.Range("L5").Value = "User cannot access restricted page"
$ws.Range("M5").Value = "High"

# The new row's steps wrap across 3 lines, so it ends up taller than a
# regular row.
$ws.Rows.Item(5).RowHeight = 63

# Column K (Test Data) got a bit wider to fit the new text.
$ws.Columns.Item(11).ColumnWidth = 21.5

# Update selection/scroll position to match where the author ended up.
$ws.Range("L5").Select()
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
